# Apply updated cryptocurrency symbol/price data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '247.58'
Set-TextValue $ws.Range('D3') '22.36'
Set-TextValue $ws.Range('D4') '5.238'
Set-TextValue $ws.Range('D5') '0.05688'
Set-TextValue $ws.Range('D6') '3.417'
Set-TextValue $ws.Range('D7') '6.313'
Set-TextValue $ws.Range('D8') '0.8068'
Set-TextValue $ws.Range('D9') '0.9011'
Set-TextValue $ws.Range('D10') '0.1422'
Set-TextValue $ws.Range('D11') '0.07433'
Set-TextValue $ws.Range('D12') '0.03054'
Set-TextValue $ws.Range('D13') '0.03075'
Set-TextValue $ws.Range('D14') '0.09383'
Set-TextValue $ws.Range('D15') '3.878'
Set-TextValue $ws.Range('D16') '0.001579'
Set-TextValue $ws.Range('D17') '0.04776'
Set-TextValue $ws.Range('D18') '0.01828'
Set-TextValue $ws.Range('D19') '0.0005802'
Set-TextValue $ws.Range('E19') '18OneONEWorstin24h'
Set-TextValue $ws.Range('D20') '0.006413'
Set-TextValue $ws.Range('D21') '0.005044'
Set-TextValue $ws.Range('D22') '0.0009969'
Set-TextValue $ws.Range('D23') '0.0001501'
Set-TextValue $ws.Range('D24') '3.689'
Set-TextValue $ws.Range('D25') '2.163'
Set-TextValue $ws.Range('D26') '0.3244'
Set-TextValue $ws.Range('D27') '0.1351'
Set-TextValue $ws.Range('D40') '0.03956'
Set-TextValue $ws.Range('B41') 'KickToken'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range('D41') '0.006793'
Set-TextValue $ws.Range('E41') '40KickTokenKICK'
Set-TextValue $ws.Range('B42') 'BKEXToken'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range('D42') '0.1066'
Set-TextValue $ws.Range('E42') '41BKEXTokenBKK'
Set-TextValue $ws.Range('B43') 'CEJI'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range('D43') '0.003201'
Set-TextValue $ws.Range('E43') '42CEJICEJI'
Set-TextValue $ws.Range('D44') '0.008756'
Set-TextValue $ws.Range('D45') '0.00005588'
Set-TextValue $ws.Range('D46') '0.00000000750'
Set-TextValue $ws.Range('D47') '0.4991'
Set-TextValue $ws.Range('D48') '0.1366'
Set-TextValue $ws.Range('D50') '0.01010'
